$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Förändrad" (Changed) date column C for all existing data
#    rows (2..221) from 45203 to 45204.
$ws.Range("C2:C221").Value = 45204

# 2. Row 221 gains an explicit row height (15, custom height) that it did
#    not have before.
$ws.Rows.Item(221).RowHeight = 15

# 3. Append a new data row (222) for case "A 47402-2023".
$ws.Range("A222").Value = "A 47402-2023"
$ws.Range("B222").Value = 45202
$ws.Range("C222").Value = 45204
$ws.Range("D222").Value = "DALARNAS LÄN"
$ws.Range("E222").Value = "AVESTA"
$ws.Range("G222").Value = 0.8
$ws.Range("H222").Value = 0
$ws.Range("I222").Value = 0
$ws.Range("J222").Value = 0
$ws.Range("K222").Value = 0
$ws.Range("L222").Value = 0
$ws.Range("M222").Value = 0
$ws.Range("N222").Value = 0
$ws.Range("O222").Value = 0
$ws.Range("P222").Value = 0
$ws.Range("Q222").Value = 0
$ws.Range("R222").Value = ""

# Match the date-format style used by the rest of column B/C, and the
# wrap-text style used by the rest of column R.
$ws.Range("B222:C222").NumberFormat = "YYYY-MM-DD"
$ws.Range("R222").WrapText = $true
